$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 58825324
$ws.Range("I12").Value = 1718.5834
$ws.Range("J12").Value = 200001980
$ws.Range("K12").Value = 1718.5834
$ws.Range("L12").Value = 200001980
$ws.Range("M12").Value = -1548.5834
$ws.Range("N12").Value = -200002320

$ws.Range("H53").Value = 1823.7646
$ws.Range("I53").Value = 2290
$ws.Range("J53").Value = 1409.3334
$ws.Range("K53").Value = 2290
$ws.Range("L53").Value = 1409.3334
$ws.Range("M53").Value = -1653
$ws.Range("N53").Value = -2683.3334

$ws.Range("H80").Value = 4231.151
$ws.Range("I80").Value = 3409.3635
$ws.Range("J80").Value = 4814.355
$ws.Range("K80").Value = 10228.0905
$ws.Range("L80").Value = 14443.065
$ws.Range("M80").Value = -9230.0905
$ws.Range("N80").Value = -16439.065

$ws.Range("H83").Value = 4231.151
$ws.Range("I83").Value = 3409.3635
$ws.Range("J83").Value = 4814.355
$ws.Range("K83").Value = 30684.2715
$ws.Range("L83").Value = 43329.19499999999
$ws.Range("M83").Value = -25692.2715
$ws.Range("N83").Value = -53313.19499999999

$ws.Range("H88").Value = 2853.6897
$ws.Range("I88").Value = 1928.1
$ws.Range("J88").Value = 3340.842
$ws.Range("K88").Value = 1928.1
$ws.Range("L88").Value = 3340.842
$ws.Range("M88").Value = -1522.1
$ws.Range("N88").Value = -4152.842000000001

$ws.Range("H91").Value = 2853.6897
$ws.Range("I91").Value = 1928.1
$ws.Range("J91").Value = 3340.842
$ws.Range("K91").Value = 1928.1
$ws.Range("L91").Value = 3340.842
$ws.Range("M91").Value = -524.0999999999999
$ws.Range("N91").Value = -6148.842000000001

$ws.Range("H132").Value = 2632.7
$ws.Range("I132").Value = 1870.04
$ws.Range("J132").Value = 6446
$ws.Range("K132").Value = 5610.12
$ws.Range("L132").Value = 19338
$ws.Range("M132").Value = -3080.12
$ws.Range("N132").Value = -24398

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3223.0938
$ws.Range("I63").Value = 2608.6875
$ws.Range("J63").Value = 3837.5
$ws.Range("K63").Value = 2608.6875
$ws.Range("L63").Value = 3837.5
$ws.Range("M63").Value = -1922.6875
$ws.Range("N63").Value = -5209.5

$ws.Range("H66").Value = 3223.0938
$ws.Range("I66").Value = 2608.6875
$ws.Range("J66").Value = 3837.5
$ws.Range("K66").Value = 13043.4375
$ws.Range("L66").Value = 19187.5
$ws.Range("M66").Value = -9611.4375
$ws.Range("N66").Value = -26051.5

$ws.Range("H110").Value = 1630.3636
$ws.Range("I110").Value = 1353.8667
$ws.Range("J110").Value = 2222.8572
$ws.Range("K110").Value = 1353.8667
$ws.Range("L110").Value = 2222.8572
$ws.Range("M110").Value = 691.1333
$ws.Range("N110").Value = -6312.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1849.9445
$ws.Range("I20").Value = 1614.5555
$ws.Range("J20").Value = 2085.3333
$ws.Range("K20").Value = 1614.5555
$ws.Range("L20").Value = 2085.3333
$ws.Range("M20").Value = -1367.5555
$ws.Range("N20").Value = -2579.3333

$ws.Range("H86").Value = 1341.86
$ws.Range("I86").Value = 1244.65
$ws.Range("J86").Value = 1730.7
$ws.Range("K86").Value = 1244.65
$ws.Range("L86").Value = 1730.7
$ws.Range("M86").Value = -121.6500000000001
$ws.Range("N86").Value = -3976.7

$ws.Range("H89").Value = 1341.86
$ws.Range("I89").Value = 1244.65
$ws.Range("J89").Value = 1730.7
$ws.Range("K89").Value = 6223.25
$ws.Range("L89").Value = 8653.5
$ws.Range("M89").Value = -607.25
$ws.Range("N89").Value = -19885.5

$ws.Range("H99").Value = 1799.9286
$ws.Range("I99").Value = 1411
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1411
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 87
$ws.Range("N99").Value = -5496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 635.5
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 671
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 671
$ws.Range("M22").Value = -250
$ws.Range("N22").Value = -1371

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3954.5454
$ws.Range("J100").Value = 3954.5454
$ws.Range("L100").Value = 11863.6362
$ws.Range("N100").Value = -13485.6362

$ws.Range("H109").Value = 2213.3333
$ws.Range("I109").Value = 1464.4
$ws.Range("J109").Value = 3149.5
$ws.Range("K109").Value = 4393.200000000001
$ws.Range("L109").Value = 9448.5
$ws.Range("M109").Value = -3353.200000000001
$ws.Range("N109").Value = -11528.5

$ws.Range("H131").Value = 1015.04877
$ws.Range("I131").Value = 827
$ws.Range("J131").Value = 1029.8948
$ws.Range("K131").Value = 2481
$ws.Range("L131").Value = 3089.6844
$ws.Range("M131").Value = 2559
$ws.Range("N131").Value = -13169.6844

$ws.Range("H132").Value = 2561.8096
$ws.Range("I132").Value = 2778.6667
$ws.Range("J132").Value = 2272.6667
$ws.Range("K132").Value = 25008.0003
$ws.Range("L132").Value = 20454.0003
$ws.Range("M132").Value = -22478.0003
$ws.Range("N132").Value = -25514.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 125810
$ws.Range("I22").Value = 1000000
$ws.Range("J22").Value = 925.7143
$ws.Range("K22").Value = 1000000
$ws.Range("L22").Value = 925.7143
$ws.Range("M22").Value = -999705
$ws.Range("N22").Value = -1515.7143

$ws.Range("H27").Value = 125810
$ws.Range("I27").Value = 1000000
$ws.Range("J27").Value = 925.7143
$ws.Range("K27").Value = 1000000
$ws.Range("L27").Value = 925.7143
$ws.Range("M27").Value = -999893
$ws.Range("N27").Value = -1139.7143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H136").Value = 27028532
$ws.Range("I136").Value = 45455664
$ws.Range("J136").Value = 2070.8667
$ws.Range("K136").Value = 136366992
$ws.Range("L136").Value = 6212.6001
$ws.Range("M136").Value = -136364442
$ws.Range("N136").Value = -11312.6001
